{"js": "// Convert M2Doc-style Word field codes (\"{ m:... }\") back into literal\n// template text runs (\"{m:...}\") - i.e. turn the complex field\n// (fldChar begin / instrText.. / fldChar end) into a plain <w:t> run\n// containing \"{\" + <trimmed field code> + \"}\".\n\nconst body = context.document.body;\nconst range = body.getRange();\nconst fields = range.fields;\nfields.load(\"items\");\nawait context.sync();\n\n// Load the instruction code for every field up front.\nfor (let i = 0; i < fields.items.length; i++) {\n  fields.items[i].load(\"code\");\n}\nawait context.sync();\n\n// Process the fields from last to first so removing one never disturbs\n// the position/identity of the fields still to be handled.\nfor (let i = fields.items.length - 1; i >= 0; i--) {\n  const f = fields.items[i];\n  let code = f.code;\n\n  // Word stores the instruction padded with exactly one leading space\n  // (right after the opening field delimiter) and one trailing space\n  // (right before the closing field delimiter). Replace those with the\n  // literal \"{\" / \"}\" template delimiters.\n  if (code.startsWith(\" \")) {\n    code = code.substring(1);\n  }\n  if (code.endsWith(\" \")) {\n    code = code.substring(0, code.length - 1);\n  }\n  const newText = \"{\" + code + \"}\";\n\n  // Identify the paragraph that owns this field before removing it.\n  f.select();\n  await context.sync();\n  const selection = context.document.getSelection();\n  const ownerParagraphs = selection.paragraphs;\n  ownerParagraphs.load(\"items\");\n  await context.sync();\n  const owner = ownerParagraphs.items[0];\n\n  // Deleting the field removes the begin/instrText/end runs while\n  // leaving the paragraph (and its pPr) untouched.\n  f.delete();\n  await context.sync();\n\n  // Re-insert the former field instruction as plain template text.\n  owner.insertText(newText, Word.InsertLocation.start);\n  await context.sync();\n}\n", "ps1": "# Convert M2Doc-style Word field codes (\"{ m:... }\") back into literal\n# template text runs (\"{m:...}\") \u2013 i.e. turn the complex field\n# (fldChar begin / instrText.. / fldChar end) into a plain <w:t> run\n# containing \"{\" + <trimmed field code> + \"}\".\n#\n# Walk the Fields collection from the end to the start so deleting a\n# field never invalidates the indices/positions of the fields still to\n# be processed.\n$d = $word.ActiveDocument\n\nfor ($i = $d.Fields.Count; $i -ge 1; $i--) {\n    $f = $d.Fields($i)\n\n    # Field.Code.Text is the raw instruction text, e.g. \" m:'Some value'.setDocumentDescription() \"\n    $code = $f.Code.Text\n\n    # Drop exactly one leading space (the one Word inserts after \"{\")\n    # and one trailing space (the one Word inserts before \"}\"), then\n    # wrap with the literal template delimiters.\n    if ($code.StartsWith(\" \")) {\n        $code = $code.Substring(1)\n    }\n    if ($code.EndsWith(\" \")) {\n        $code = $code.Substring(0, $code.Length - 1)\n    }\n    $newText = \"{\" + $code + \"}\"\n\n    # Remember where the field starts so we can insert the replacement\n    # text at exactly that spot.\n    $insertionPoint = $f.Code.Start - 1\n\n    # Removing the field removes the begin/instrText/end runs and\n    # leaves the (now empty) paragraph / surrounding runs untouched.\n    $f.Delete()\n\n    $r = $d.Range($insertionPoint, $insertionPoint)\n    $r.InsertAfter($newText)\n}\n\nWrite-Output \"done\"\n"}
